$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Fix up sheet1 view: deselect tab, clear topLeftCell, set column D selected
$ws1.Activate()
$ws1.Range("D1:D1048576").Select()

# Add the new worksheet after Hárok1
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Hárok2"

# Column widths (closest achievable via ColumnWidth rounding)
$newSheet.Columns.Item(1).ColumnWidth = 16.592447916666668
$newSheet.Columns.Item(2).ColumnWidth = 12.307291666666666

# Header row
$newSheet.Range("A1").Value = "no-analyzers"
$newSheet.Range("B1").Value = "bh-analyzers"

# Data rows 2-102
$aVals = @(5.2077199183333338,5.0867765716666664,5.124200506666666,5.0411099300000002,5.0984446649999997,5.0911279066666673,5.03698259,5.1205577266666671,5.0111799749999992,5.144822445,5.0362942083333335,5.0457826033333335,5.0937298283333341,5.0533582233333334,5.0754201183333336,5.0227803583333337,5.1038436416666659,5.1383525600000004,5.1521923049999998,5.0903149933333331,5.1339030983333327,5.0846932200000001,5.1300538500000004,5.1126851116666669,5.1000856599999995,5.0916144683333338,5.0593575116666667,5.0800218216666666,6.4831331399999996,5.1219881783333339,5.047798573333333,5.0676044849999995,5.0778693733333329,5.0700108966666662,5.0942955716666667,5.1340936033333335,5.0784151616666673,5.1243361416666673,5.0218895650000004,5.0797974083333335,5.1266723033333337,5.0952018899999993,5.0937159766666671,5.0493276666666667,5.1136617700000002,5.0776944833333335,5.1114495866666667,5.1029539633333334,5.1071097100000005,5.1154488783333338,5.1339069733333336,5.0680982999999999,5.1170470499999992,5.0816471583333334,5.0454565666666662,5.0343351316666665,5.079235726666667,5.1073750250000005,5.0580783133333336,5.0529728,5.0188409366666669,5.1006888799999999,5.0318521500000006,5.1227026483333331,5.0512788149999999,5.1546176216666666,5.0919838716666668,5.1133705016666662,5.1071462866666666,5.0892534316666671,5.1063405933333339,5.1199859266666667,5.0507365066666674,5.115553218333333,5.136543913333333,5.1064642583333333,5.1213796466666661,5.1012666366666659,5.075045358333333,5.0493316066666667,5.0889768499999999,5.0795654483333337,5.1182779016666666,5.1066178249999998,5.0434238250000005,5.1141631566666668,5.02961641,5.1255421549999998,5.1252731933333333,5.053311935,5.1512297366666662,5.0790847599999998,5.1026256783333332,5.1065256683333331,5.1178963300000007,5.0448302783333334,5.1006072316666664,5.0917179683333327,5.0835318916666674,5.0680739033333335,5.0915999686842124)
$bVals = @(5.3805436433333327,5.2604513500000003,5.2444562966666668,5.206638645,5.2595761383333333,5.3672456433333329,5.2267668299999999,5.267906393333333,5.24258778,5.2803687816666658,5.2803269166666666,5.2261645000000003,5.2651057016666671,5.2709320450000003,5.2227375450000002,5.2372609016666667,5.2434614083333333,5.1849784866666662,5.2061898166666669,5.2519269550000001,5.2707379633333327,5.2505220783333337,5.1730502533333338,5.1882614550000001,5.1716625816666664,5.1836574200000003,5.2357000249999999,5.1881177666666662,5.225218176666667,5.1867903783333329,5.1956955983333337,5.2541469233333338,5.1797570333333338,5.2868905500000007,5.2219127066666671,5.2204412699999994,5.2627663583333328,5.1819683183333334,5.2254005883333328,5.2265904983333336,5.2560468266666662,5.1975453633333331,5.1796197133333335,5.227595851666667,5.224930866666667,5.2042806666666674,5.1855446399999998,5.1944351733333338,5.2824494499999997,5.2305388033333333,5.2484029700000008,5.2499754983333338,5.2206572283333337,5.2516650699999996,5.2035187249999995,5.1800755299999999,5.2635697166666668,5.2527042133333337,5.2280753416666661,5.261459181666666,5.2046220933333327,5.2023383699999997,5.2597090466666669,5.2199192883333341,5.2683315016666672,5.2810217899999996,5.2094744266666666,5.2086711183333332,5.2919696483333336,5.2841255949999999,5.2771404500000001,5.2510063249999996,5.2533049900000002,5.2019445683333334,5.230328861666667,5.1693880549999998,5.2453743866666667,5.2884491833333334,5.2157062283333335,5.2034156399999993,5.2583775433333333,5.2792097650000001,5.2632653600000001,5.2715082083333336,5.1916642666666668,5.2063926233333335,5.210658051666667,5.2069800866666665,5.2550031666666666,5.2552174450000004,5.2190040283333339,5.2500028233333333,5.2436222716666672,5.2306272266666669,5.2483529433333338,5.2581898049999998,5.3075590199999994,5.2175166150000001,5.2362783766666672,5.2204477216666669,5.2344651750438596)
for ($i = 0; $i -lt $aVals.Length; $i++) {
    $r = $i + 2
    $newSheet.Cells.Item($r, 1).Value = $aVals[$i]
    $newSheet.Cells.Item($r, 2).Value = $bVals[$i]
}

# Summary rows 104-105 (column B only)
$newSheet.Range("B104").Value = 1.0280590005574548
$newSheet.Range("B105").Value = 8.5719123815788301

# Final selection / active sheet
$newSheet.Activate()
$newSheet.Range("F8").Select()
